$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the table name used in the "Insert Game" section's INSERT statement:
#    "INSERT INTO games (...)" -> "INSERT INTO game (...)"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "INSERT INTO games (game_id, game_name, game_img, description, price, quantity, create_date) VALUES",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Replace only the word "games" (it sits right after "INSERT INTO ") so the
    # rest of the paragraph/run is left untouched.
    $gamesStart = $rng.Start + "INSERT INTO ".Length
    $gamesEnd = $gamesStart + "games".Length
    $gamesRange = $d.Range($gamesStart, $gamesEnd)
    $gamesRange.Text = "game"
}

# ---------------------------------------------------------------------------
# 2) Insert 10 new rows (games 23-32) right after the INSERT INTO statement,
#    before the existing "(10, 'The Legend of Zelda...')" paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "INSERT INTO game (game_id, game_name, game_img, description, price, quantity, create_date) VALUES",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $lines = @(
        "(23, 'Stardew Valley', 'game23.jpg', 'Một trò chơi mô phỏng cuộc sống nông trại với các hoạt động như trồng trọt, câu cá, và kết bạn.', 150000, 120, '2021-07-10'),",
        "(24, 'Cyberpunk 2077', 'game24.jpg', 'Một trò chơi nhập vai hành động thế giới mở, diễn ra trong thành phố tương lai Night City.', 320000, 95, '2020-12-15'),",
        "(25, 'Among Us', 'game9.jpg', 'Một trò chơi nhiều người chơi, nơi bạn phải tìm ra kẻ phản bội trong nhóm.', 50000, 300, '2023-04-20'),",
        "(26, 'Elden Ring', 'game10.jpg', 'Một trò chơi nhập vai hành động thế giới mở với sự hợp tác giữa FromSoftware và George R.R. Martin.', 350000, 80, '2022-02-25'),",
        "(27, 'Terraria', 'game1.jpg', 'Một trò chơi sandbox phiêu lưu hành động 2D với khả năng xây dựng và khám phá.', 100000, 200, '2020-05-16'),",
        "(28, 'Genshin Impact', 'game2.jpg', 'Một trò chơi nhập vai hành động thế giới mở với cơ chế gacha.', 0, 500, '2023-10-12'),",
        "(29, 'League of Legends', 'game3.jpg', 'Một trò chơi MOBA nổi tiếng với các trận đấu đội cạnh tranh.', 0, 1000, '2019-09-23'),",
        "(30, 'Apex Legends', 'game4.jpg', 'Một trò chơi battle royale với các nhân vật có kỹ năng độc đáo.', 0, 400, '2021-11-18'),",
        "(31, 'Fall Guys', 'game5.jpg', 'Một trò chơi nhiều người chơi với các vòng đấu vượt chướng ngại vật vui nhộn.', 80000, 350, '2022-03-05'),",
        "(32, 'The Sims 4', 'game6.jpg', 'Một trò chơi mô phỏng cuộc sống, nơi bạn tạo và quản lý các nhân vật Sims của mình.', 180000, 150, '2018-10-20'),"
    )

    $insertAfter = $rng.End
    foreach ($line in $lines) {
        $insertPoint = $d.Range($insertAfter, $insertAfter)
        $insertPoint.InsertParagraphAfter()
        $newParaStart = $insertAfter + 1
        $newParaRange = $d.Range($newParaStart, $newParaStart)
        $newParaRange.Text = $line
        $insertAfter = $newParaStart + $line.Length
    }
}
